$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert a new column A (TabName) and a new row 3 (FilesTab) ---
$ws.Columns("A").Insert()
$ws.Rows(3).Insert()

# --- Set cell values in the same order the original author typed them, so
#     the shared-string table is built up in the same sequence:
#     TabName, CasesTab, FilesTab, the two count queries, then the two
#     full queries. ---

# Row 1 headers
$ws.Range("A1").Value2 = "TabName"
$ws.Range("B1").Value2 = "query"
$ws.Range("C1").Value2 = "StatQuery"
$ws.Range("D1").Value2 = "dbExcel"
$ws.Range("E1").Value2 = "WebExcel"

# Tab name labels first
$ws.Range("A2").Value2 = "CasesTab"
$ws.Range("A3").Value2 = "FilesTab"

# Count queries (StatQuery column) next
$casesCountQuery = @'
MATCH (f:file)
OPTIONAL MATCH (f)-[*]->(a:arm)-[:of_trial]->(ct:clinical_trial)
OPTIONAL MATCH (f)-[*]->(c:case)
WITH f,a,ct,c
      WHERE f.file_type = 'Aligned DNA reads file'
RETURN
    COUNT(DISTINCT ct.clinical_trial_designation) AS Trials,
    COUNT(DISTINCT c.case_id) AS Cases,
    COUNT(DISTINCT f) AS Files
'@
$ws.Range("C2").Value2 = $casesCountQuery

$filesCountQuery = @'
MATCH (f:file)
OPTIONAL MATCH (f)-[*]->(a:arm)-[:of_trial]->(ct:clinical_trial)
OPTIONAL MATCH (f)-[*]->(c:case)
WITH f,a,ct,c
        WHERE f.file_type = 'Aligned DNA reads file'
RETURN
    COUNT(DISTINCT ct.clinical_trial_designation) AS Trials,
    COUNT(DISTINCT c.case_id) AS Cases,
    COUNT(DISTINCT f) AS Files
'@
$ws.Range("C3").Value2 = $filesCountQuery

# Full queries last
$casesQuery = @'
MATCH (c:case)
 MATCH (c)-[:of_arm]->(a:arm)-[:of_trial]->(ct:clinical_trial)
 MATCH (f:file)-[*]->(c)
WHERE f.file_type = 'Aligned DNA reads file' 
RETURN DISTINCT
    c.case_id AS `Case ID`,
     ct.clinical_trial_designation AS `Trial Code`,
     a.arm_id AS Arm,
      a.arm_drug AS `Arm Treatment`,
c.disease AS Diagnosis,
  c.gender AS Gender,
    c.race AS Race,
    c.ethnicity AS Ethnicity
'@
$ws.Range("B2").Value2 = $casesQuery

$filesQuery = @'
MATCH (f:file)
OPTIONAL MATCH (f)-[*]->(a:arm)-[:of_trial]->(ct:clinical_trial)
OPTIONAL MATCH (f)-[*]->(c:case)
OPTIONAL MATCH (f)-->(parent)
WITH f,a,ct,c,parent
WHERE f.file_type = 'Aligned DNA reads file'
 WITH
    f, parent, c, a, ct,
    ['Bytes', 'KB', 'MB', 'GB', 'TB'] AS units,
    toInteger(floor(log(f.file_size)/log(1024))) as i,
    2 as precision
WITH
    f, parent, c, a, ct,
    f.file_size /(1024^i) AS value,
    10^precision AS factor,
    units[i] as unit
WITH
    f, parent, c, a, ct, unit,
    round(factor * value)/factor AS size
RETURN DISTINCT
    f.file_name AS `File Name`,
    head(labels(parent)) as Association,
    f.file_description AS Description,
    f.file_format AS `File Format`,
    CASE size % 1 WHEN 0 THEN apoc.convert.toInteger(size)+' ' +unit ELSE size+' ' +unit END AS Size,
    ct.clinical_trial_designation AS `Trial Code`,
    a.arm_id AS Arm,
    c.case_id AS `Case ID`
'@
$ws.Range("B3").Value2 = $filesQuery

# File-name columns (reuse existing shared strings)
$ws.Range("D2").Value2 = "TC01_Trials_Filter_AssocFileType-AlignedDNA_Neo4jData.xlsx"
$ws.Range("E2").Value2 = "TC01_Trials_Filter_AssocFileType-AlignedDNA_WebData.xlsx"
$ws.Range("D3").Value2 = "TC01_Trials_Filter_AssocFileType-AlignedDNA_Neo4jData.xlsx"
$ws.Range("E3").Value2 = "TC01_Trials_Filter_AssocFileType-AlignedDNA_WebData.xlsx"

# --- Wrap text on the long query cells (matches style index 1 / "Normal 2") ---
$ws.Range("B2").WrapText = $true
$ws.Range("C2").WrapText = $true
$ws.Range("B3").WrapText = $true
$ws.Range("C3").WrapText = $true

# --- Column widths (B-E already carry the correct original widths after the
#     column insert shifted them; only the new column A needs to be set) ---
$ws.Columns("A").ColumnWidth = 8

# --- Row heights ---
$ws.Rows(2).RowHeight = 188.5
$ws.Rows(3).RowHeight = 409.5

# --- View: zoom + selection ---
$excel.ActiveWindow.Zoom = 85
$ws.Range("C2").Select()
